$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet tracked "Monivalinta" (multiple-choice) option groups on rows
# 3-6. Updating it "for year 2018" replaces the content with a single,
# current option group (the "Toteuma" / realisation group, columns A-F)
# and removes the option groups that used to live in rows 3-6.

$ws.Range("A3:H6").ClearContents()

# Touch row height on rows 3-6 so the now-empty rows 4-6 are still kept
# around (instead of being dropped entirely) once their content is gone.
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(3).RowHeight
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(4).RowHeight
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(5).RowHeight
$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(6).RowHeight

# Row 3 now holds the "Toteuma" (realisation) option group.
$ws.Range("A3").Value = "Monivalinta"
$ws.Range("B3").Value = "Toteuma"
$ws.Range("C3").Value = "pkv"
$ws.Range("D3").Value = "Ei toteutunut"
$ws.Range("E3").Value = "Osittain toteutunut"
$ws.Range("F3").Value = "Toteutunut"

# Leave the selection where the edit took place.
$ws.Range("C3").Select()
